# Update the "ManageContactPage" worksheet with the new test data values
# (method names were renamed in the corresponding Java page-object class,
# so the sample/expected values used by the tests were refreshed too).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ManageContactPage")

$ws.Range("A2").Value = 8281828182
$ws.Range("A3").Value = "userdummy@yopmail.com"
$ws.Range("A4").Value = "Flat No 1207A"
$ws.Range("A5").Value = 150
$ws.Range("A6").Value = 60

# Widen column A to fit the new, longer values.
$ws.Columns.Item(1).ColumnWidth = 19.42

# Move the active selection to D6, matching the final saved state.
$ws.Activate()
$ws.Range("D6").Select()
